# Commit: "Changing maps a bit"
# The Feuil1 sprite/map grid had several cells (columns G, I, K, T, V, X
# across rows 8-23) marked with sprite ids 1 / 4. This edit clears them
# back to 0 ("changing the maps a bit"), then leaves the view zoomed to
# 100% with V22 as the new active selection (previously X23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Feuil1" is the tabSelected / active sheet

# Columns that go from 1 -> 0 for every row 8..23
$colsToClear = @("G", "K", "T", "X")
for ($row = 8; $row -le 23; $row++) {
    foreach ($col in $colsToClear) {
        $ws.Range("$col$row").Value = 0
    }
}

# Rows 9 and 22 also had I/V marked with sprite id 4 -> clear those too
foreach ($row in @(9, 22)) {
    $ws.Range("I$row").Value = 0
    $ws.Range("V$row").Value = 0
}

# Update the view: zoom to 100% and move the selection to V22
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("V22").Select()
